$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting rows 41-57 down to 42-58.
$ws.Rows.Item(41).Insert()

$ws.Range("A41").Value = 5
$ws.Range("B41").Value = "Macroferia Regional de Talca"
$ws.Range("C41").Value = "Maule"
$ws.Range("D41").Value = 44813
$ws.Range("D41").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E41").Value = 7
$ws.Range("F41").Value = 100112040
$ws.Range("G41").Value = "Cilantro"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 150
$ws.Range("K41").Value = 7500
$ws.Range("L41").Value = 7500
$ws.Range("M41").Value = 7500
$ws.Range("N41").Value = "$/caja 36 atados"
$ws.Range("O41").Value = "Región del Maule"
$ws.Range("P41").Value = 208
$ws.Range("Q41").Value = 36
$ws.Range("R41").Value = "Hortaliza"
